# Sound attenuating chamber small BOM - update LED driver / resistor parts,
# and remove the obsolete Doric_FRJ_rotary_mount row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 29: IR LED module -> more specific part description (part number unchanged)
$ws.Range("B29").Value = "IR LED module, 730nm, 9.2V @ 1A"

# Row 30: 12 Ohm resistor 7W -> LED driver (with new part number)
$ws.Range("B30").Value = "LED driver, 12V input 700mA CC output"
$ws.Range("E30").Value = "176-6768"

# Remove the last two rows (header + Doric_FRJ_rotary_mount / 5mm clear Acrylic entry)
$ws.Range("B40:D41").EntireRow.Delete()

# Update the active cell selection to match the saved view state
$ws.Activate()
$ws.Range("G25").Select()

$wb.Save()
